# Drop the header row (phone/nama) now that the importer indexes columns
# directly instead of relying on a header label row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Delete()

# Reselect the new top row (now the first data row) the way Excel leaves
# the selection after a row deletion/header removal.
$ws.Rows(1).Select()

# Outline max level drops now that one fewer row level exists.
$ws.Outline.ShowLevels(1)
